# ---------------------------------------------------------------------------
# Commit: "20DIC2022: ejemplos Algoritmo Regresion Logistica"
#
# Restructures the three "EDA conclusions" bullet paragraphs right after
# "De acuerdo al análisis EDA inicial se determina lo siguiente:":
#   1) the "Todos los atributos son categóricos..." bullet is replaced by the
#      "El atributo unacc..." text (which used to be bullet #2), with two new
#      clauses added about the class attribute counts,
#   2) the old bullet #2 ("El atributo unacc...") becomes a short wrap-up
#      sentence about the dataset being an unbalanced multiclass
#      classification problem (still carries the original list pPr, per the
#      source XML),
#   3) the first of the three trailing empty/bold paragraphs becomes a new
#      "Finalmente se comenta..." bullet (an expanded version of the old
#      "Todos los atributos..." bullet), followed by fresh narrative
#      paragraphs about the train/test split and a lead-in to the upcoming
#      algorithm examples.
#
# Paragraph identities are resolved by distinctive text prefix *before* any
# edit is applied (so a later replacement's own new text can never be
# mistaken for an earlier anchor), then edits are applied in document order
# using the captured Paragraph objects.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

function Find-ParagraphByPrefix($doc, [string]$prefix, [int]$fallbackIndex) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $candidate = $doc.Paragraphs.Item($i)
        if ($candidate.Range.Text.StartsWith($prefix)) {
            return $candidate
        }
    }
    return $doc.Paragraphs.Item($fallbackIndex)
}

# Resolve all three anchors up front, against the untouched document, so
# that editing paragraph 1's text cannot shadow the lookup for paragraph 2.
$p1 = Find-ParagraphByPrefix $d "Todos los atributos son categ" 85
$p2 = Find-ParagraphByPrefix $d "El atributo" 86
$p3 = $d.Paragraphs.Item(87)

# --- Bullet 1: "Todos los atributos son categóricos..." -> "El atributo ... unacc..." ---
# Keeps its original list pPr (pStyle Prrafodelista / numPr ilvl0,numId7 / ind left426).
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:ind w:left="426"/></w:pPr><w:r><w:t>El atributo “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>unacc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>” (inaceptable),</w:t></w:r><w:r><w:t xml:space="preserve"> gestiona la cantidad de valores </w:t></w:r><w:r><w:t>más</w:t></w:r><w:r><w:t xml:space="preserve"> alto en comparación con los otros atributos</w:t></w:r><w:r><w:t xml:space="preserve"> que gestiona a clase</w:t></w:r><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>a</w:t></w:r><w:r><w:t>cc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>g</w:t></w:r><w:r><w:t>ood</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>vgood</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r><w:r><w:t>, lo cual nos lle</w:t></w:r><w:r><w:t>va a determinar que este</w:t></w:r><w:r><w:t xml:space="preserve"> atributo del</w:t></w:r><w:r><w:t xml:space="preserve"> dataset esta desbalanceado a comparación </w:t></w:r><w:r><w:t>a los otros atributos de la clase.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p1.Range.InsertXML($xml1)

# --- Bullet 2: "El atributo ... unacc..." -> "Con esta información..." ---
# Also keeps its original list pPr (unchanged in the source diff).
$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:ind w:left="426"/></w:pPr><w:r><w:t xml:space="preserve">Con esta información se evidencia que </w:t></w:r><w:r><w:t>este es un problema de clasificación multiclase desequilibrado.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p2.Range.InsertXML($xml2)

# --- First trailing empty/bold paragraph becomes the new bullet 3 plus the
#     freshly-authored narrative paragraphs that follow it. Word's InsertXML
#     folds the very last paragraph mark of the inserted fragment into the
#     host boundary, so a trailing run of N empty paragraphs needs N+1
#     <w:p/> markers in the payload to materialize as N paragraphs. ---
$xml3 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:ind w:left="426"/></w:pPr><w:r><w:t>Finalmente se comenta que t</w:t></w:r><w:r><w:t>odos los atributos son</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>categóric</w:t></w:r><w:r><w:t>o</w:t></w:r><w:r><w:t xml:space="preserve">s, por lo que para poder ejecutar los algoritmos de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sklearn</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> debemos convertir estos datos categóricos en valores enteros.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Una vez ejecutado el </w:t></w:r><w:r><w:t xml:space="preserve">proceso de transformación de </w:t></w:r><w:r><w:t>valores categóricos a numéricos</w:t></w:r><w:r><w:t xml:space="preserve">, procedemos a </w:t></w:r><w:r><w:t xml:space="preserve">dividir el dataset resultante en datos de entrenamiento </w:t></w:r><w:r><w:t>y datos de test</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>una práctica común es 80% de los dat</w:t></w:r><w:r><w:t xml:space="preserve">os como entrenamiento y el 20% restante como </w:t></w:r><w:r><w:t>datos de test</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> e</w:t></w:r><w:r><w:t xml:space="preserve">sto evitará problemas en los que nuestro algoritmo pueda fallar por </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sobregenerar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> conocimiento</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Una vez generados estos datos se procede a ejecutar </w:t></w:r><w:r><w:t xml:space="preserve">los algoritmos </w:t></w:r></w:p><w:p/><w:p/><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p3.Range.InsertXML($xml3)
